# Generate Report for Handback
#
# The "ed31471e-b328-441f-a105-ede1d361a2df.md" file has now been handed
# back (it was previously only "Ready for handoff"). Update its status on
# every sheet to match the already-handed-back file, and stamp the latest
# handback datetime for each locale.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 is the ed31471e-...-.md file ---
$overview.Range("B3").Value = "Handed back: in sync with en-US"
$overview.Range("C3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the ed31471e-...-.md file ---
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("H3").Value = "2016-03-21 20:47:28"

# --- de-de sheet: row 3 is the ed31471e-...-.md file ---
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("H3").Value = "2016-03-21 20:47:34"
